$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(376).Insert()

$ws.Range("A376").Value = 1
$ws.Range("B376").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C376").Value = 'Arica y Parinacota'
$ws.Range("D376").Value = 45106
$ws.Range("E376").Value = 15
$ws.Range("F376").Value = 'Fruta'
$ws.Range("G376").Value = 100108
$ws.Range("H376").Value = 'Tropicales y subtropicales'
$ws.Range("I376").Value = 100108006
$ws.Range("J376").Value = 'Plátano'
$ws.Range("K376").Value = 'Sin especificar'
$ws.Range("L376").Value = 'Pintón'
$ws.Range("M376").Value = 120
$ws.Range("N376").Value = 15000
$ws.Range("O376").Value = 16000
$ws.Range("P376").Value = 15500
$ws.Range("Q376").Value = '$/caja 20 kilos'
$ws.Range("R376").Value = 'Ecuador'
$ws.Range("S376").Value = 775
$ws.Range("T376").Value = 20
